# Revise budget for 0.75 in boards and add cut list.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Item text updates first, in the order the new strings were authored ---
# (European beech boards go from 1/2" to 3/4" stock)
$ws.Range("A3").Value = "European beech 3/4"" boards"
# (5/8" plywood panel swapped out for a 3/4" baltic birch panel)
$ws.Range("A2").Value = "3/4"" baltic birch plywood 30""x30"" panel"
# (new 3/4" natural birch plywood panel item)
$ws.Range("A4").Value = "3/4"" natural birch plywood 2'x4' panel"
# (cut list retitled for beech)
$ws.Range("H1").Value = "Beech cut list"

# --- Row 2 (qty/price update) ---
$ws.Range("B2").Value = "panel"
$ws.Range("D2").Value = 28

# --- Row 3 (cut list thickness 0.5 -> 0.75, quantity now a literal, plus board footage formula) ---
$ws.Range("B3").Value = "board foot"
$ws.Range("C3").Value = 8
$ws.Range("H3").Value = 0.75
$ws.Range("O3").Formula = "=SUM(N3:N7)/144"

# --- Row 4 (Lumber delivery moves to row 5; new plywood panel item moves in) ---
$ws.Range("B4").Value = "panel"
$ws.Range("D4").Value = 29
$ws.Range("H4").Value = 0.75

# --- Row 5 (now holds what used to be the "Lumber delivery" budget line) ---
$ws.Range("A5").Value = "Lumber delivery"
$ws.Range("C5").Value = 1
$ws.Range("D5").Value = 40
$ws.Range("E5").Formula = "=C5*D5"
$ws.Range("H5").Value = 0.75

# --- Row 6 (budget D/E columns no longer used here; cut list updates) ---
$ws.Range("D6").Clear()
$ws.Range("E6").Clear()
$ws.Range("H6").Value = 0.75
$ws.Range("I6").Value = 4

# --- Row 7 (adds a budget subtotal; cut list thickness update) ---
$ws.Range("E7").Formula = "=SUM(E2:E5)"
$ws.Range("H7").Value = 0.75

# --- Selection / view state ---
$ws.Range("G1:M7").Select() | Out-Null
